$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.295.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.382.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.383.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000168"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.812.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.969.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.374.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +11.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.55"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.07"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.85"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.27"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "566.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.496.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.20"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.67"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.41"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.588"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.01%  "
